$d = $word.ActiveDocument

# 1. Replace the subtitle text "RequestSolved!" with the new title text.
$d.Content.Find.Execute("RequestSolved!", $true, $false, $false, $false, $false, $true, 1, $false,
                         "E-commerce de joias e itens de artesanato", 2) | Out-Null

# 2. Remove the orange "Tema" table that used to follow the subtitle.
$d.Tables.Item(1).Delete()

# 3. The table left behind an empty paragraph that is no longer needed (the
#    paragraph mark right after the now-renamed subtitle paragraph); merge it
#    away so the subtitle paragraph is immediately followed by the body text.
$full = $d.Content.Text
$idx = $full.IndexOf("E-commerce de joias e itens de artesanato")
$pos = $idx + ("E-commerce de joias e itens de artesanato").Length + 1
$d.Range($pos, $pos + 1).Delete() | Out-Null

# 4. Add commas around "por exemplo" and fix the "arvores" -> "árvores" typo.
$d.Content.Find.Execute(
    "como por exemplo vasos para plantas, arvores feitas com fios e pedras",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "como, por exemplo, vasos para plantas, árvores feitas com fios e pedras", 2) | Out-Null

# 5. Add commas around "por exemplo" and drop the colon in the other sentence.
$d.Content.Find.Execute(
    "artesanais como por exemplo: festa das cerejeiras",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "artesanais como, por exemplo, festa das cerejeiras", 2) | Out-Null

# 6. Extend the sentence about social media with the new closing clause.
$d.Content.Find.Execute(
    "comunicação e compartilhamento de histórias.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "comunicação e compartilhamento de histórias com o intuito de atingir mais pessoas com os produtos produzidos.", 2) | Out-Null

# 7. The stale rendered-page-break marker in front of "Troca;" is dropped by
#    touching the run (Word only keeps that cache hint while the run is
#    untouched).
$d.Content.Find.Execute("Troca;", $true, $false, $false, $false, $false, $true, 1, $false,
                         "Troca;", 2) | Out-Null

# 8. Fix the "fara" -> "fará" typo.
$d.Content.Find.Execute(
    "que fara o cálculo do frete.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "que fará o cálculo do frete.", 2) | Out-Null
